$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1700113000516159
$ws.Range("H2").Value = 79.48660622907647
$ws.Range("I2").Value = -9.606803240865194
$ws.Range("G3").Value = 0.09553884385491446
$ws.Range("H3").Value = 44.55015471770156
$ws.Range("G4").Value = 0.01041703237206153
$ws.Range("H4").Value = 11.11184087113988
$ws.Range("G5").Value = 0.03018724122957809
$ws.Range("H5").Value = 325.4493226611183
$ws.Range("G6").Value = -0.2073083101560849
$ws.Range("H6").Value = 6.278322139331491
$ws.Range("G7").Value = -0.1720103004699615
$ws.Range("H7").Value = 31.1616518203774
$ws.Range("G8").Value = -0.3967428515732616
$ws.Range("H8").Value = -7.153696669675853
$ws.Range("G9").Value = -0.3903293101199335
$ws.Range("H9").Value = 2.083377697258999
$ws.Range("G10").Value = -0.03938610559419427
$ws.Range("H10").Value = -343.0507692776858
$ws.Range("G11").Value = 0.131229521380489
$ws.Range("H11").Value = 916.8218913103194
$ws.Range("G12").Value = 0.2362752755454559
$ws.Range("H12").Value = 3.999525695716288
$ws.Range("G13").Value = 0.248005697511277
$ws.Range("H13").Value = -5.823214628183954
$ws.Range("G14").Value = -0.0384776241020057
$ws.Range("H14").Value = -302.1796655296345
$ws.Range("G15").Value = 0.02573307032820443
$ws.Range("H15").Value = 27.47526208952876
$ws.Range("G16").Value = 0.1647343929924573
$ws.Range("H16").Value = 39.60647623517845
$ws.Range("G17").Value = 0.2150641022976667
$ws.Range("H17").Value = -1.722236346511925
$ws.Range("G18").Value = 0.05223190306592695
$ws.Range("H18").Value = -13.61859264406413
$ws.Range("G19").Value = 0.06848999353972837
$ws.Range("H19").Value = -23.97544390207544
$ws.Range("G20").Value = -0.1788752086667082
$ws.Range("H20").Value = -22.91212145254338
$ws.Range("G21").Value = -0.1838935699811558
$ws.Range("H21").Value = 7.974944840722522
$ws.Range("G22").Value = 0.04237595348718298
$ws.Range("H22").Value = -22.08618723903151
$ws.Range("G23").Value = 0.05559011495589583
$ws.Range("H23").Value = 36.11634996571488
$ws.Range("G24").Value = 0.1307625426463148
$ws.Range("H24").Value = 12.98510174117692
$ws.Range("G25").Value = 0.1622800922161463
$ws.Range("H25").Value = 6.709843581613406
$ws.Range("G26").Value = 0.005713554100652987
$ws.Range("H26").Value = -89.19458906613869
$ws.Range("G27").Value = 0.03574884830703168
$ws.Range("H27").Value = -29.16364222290547
$ws.Range("G28").Value = 0.1446675415400978
$ws.Range("H28").Value = -5.391062129033571
$ws.Range("G29").Value = 0.1504233297723364
$ws.Range("H29").Value = -11.88013550758164
$ws.Range("G30").Value = 0.03686650458446161
$ws.Range("H30").Value = 88.42042774653645
$ws.Range("G31").Value = 0.03731925139650764
$ws.Range("H31").Value = 284.5375006907184
$ws.Range("G32").Value = -0.02204040112013462
$ws.Range("H32").Value = -159.1002069863162
$ws.Range("G33").Value = 0.001302360005632001
$ws.Range("H33").Value = -95.0108105155411
$ws.Range("G34").Value = 0.09656640136724229
$ws.Range("H34").Value = -24.53863867869114
$ws.Range("G35").Value = 0.1697810557233021
$ws.Range("H35").Value = 31.96041747957758
$ws.Range("G36").Value = -0.04102729632469536
$ws.Range("H36").Value = -372.9350949672348
$ws.Range("G37").Value = -0.02215038138807492
$ws.Range("H37").Value = -244.6367701321524
$ws.Range("G38").Value = -0.05827896936535189
$ws.Range("H38").Value = -2753.436233478201
$ws.Range("G39").Value = -0.02541269692003024
$ws.Range("H39").Value = 23.93506133103244
$ws.Range("G40").Value = 0.1364616655626349
$ws.Range("H40").Value = -7.514335647230055
$ws.Range("G41").Value = 0.1350205909486641
$ws.Range("H41").Value = -16.3437081227754
$ws.Range("G42").Value = 0.07116193170672798
$ws.Range("H42").Value = 10.21766792845759
$ws.Range("G43").Value = 0.06067289311419088
$ws.Range("H43").Value = 74.54539673314883
$ws.Range("G44").Value = 0.04280101998521783
$ws.Range("H44").Value = 203.2783894440889
$ws.Range("G45").Value = 0.03808508437875471
$ws.Range("H45").Value = -7.240857612930888
$ws.Range("G46").Value = -0.06522490415401151
$ws.Range("H46").Value = 0.9039214339372388
$ws.Range("G47").Value = -0.0939668769660478
$ws.Range("H47").Value = -127.4692640543958
$ws.Range("G48").Value = -0.0989509775500266
$ws.Range("H48").Value = 21.45214847615965
$ws.Range("G49").Value = -0.1107892339138908
$ws.Range("H49").Value = 43.89890378583468
$ws.Range("G50").Value = 0.04625888833092191
$ws.Range("H50").Value = -57.51140034684327
$ws.Range("G51").Value = 0.1201541862885204
$ws.Range("H51").Value = 19.82965363615456
$ws.Range("G52").Value = 0.04571569971152875
$ws.Range("H52").Value = -23.32031925198192
$ws.Range("G53").Value = 0.05906170640851614
$ws.Range("H53").Value = -12.55646648909768
$ws.Range("G54").Value = -0.06740290996951226
$ws.Range("H54").Value = 3.600757242623081
$ws.Range("G55").Value = -0.04325139647440099
$ws.Range("H55").Value = 44.00113543441045
$ws.Range("G56").Value = 0.06389822096236623
$ws.Range("H56").Value = 39.43638620802007
$ws.Range("G57").Value = 0.09823300568525047
$ws.Range("H57").Value = 1799.956042243333
